$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before current row 215 (old rows 215-314 shift down to 217-316).
$ws.Rows.Item(215).Resize(2).Insert()

# --- New row 215 ---
$ws.Cells.Item(215, 1).Value = 11
$ws.Cells.Item(215, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(215, 3).Value = "Bíobío"
$ws.Cells.Item(215, 4).Value = 44523
$ws.Cells.Item(215, 5).Value = 8
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100102
$ws.Cells.Item(215, 8).Value = "Cítricos"
$ws.Cells.Item(215, 9).Value = 100102003
$ws.Cells.Item(215, 10).Value = "Limón"
$ws.Cells.Item(215, 11).Value = "Sin especificar"
$ws.Cells.Item(215, 12).Value = "1a amarillo"
$ws.Cells.Item(215, 13).Value = 450
$ws.Cells.Item(215, 14).Value = 5500
$ws.Cells.Item(215, 15).Value = 6000
$ws.Cells.Item(215, 16).Value = 5778
$ws.Cells.Item(215, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(215, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(215, 19).Value = 321
$ws.Cells.Item(215, 20).Value = 18

# --- New row 216 ---
$ws.Cells.Item(216, 1).Value = 11
$ws.Cells.Item(216, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(216, 3).Value = "Bíobío"
$ws.Cells.Item(216, 4).Value = 44523
$ws.Cells.Item(216, 5).Value = 8
$ws.Cells.Item(216, 6).Value = "Fruta"
$ws.Cells.Item(216, 7).Value = 100102
$ws.Cells.Item(216, 8).Value = "Cítricos"
$ws.Cells.Item(216, 9).Value = 100102003
$ws.Cells.Item(216, 10).Value = "Limón"
$ws.Cells.Item(216, 11).Value = "Sin especificar"
$ws.Cells.Item(216, 12).Value = "2a amarillo"
$ws.Cells.Item(216, 13).Value = 350
$ws.Cells.Item(216, 14).Value = 5500
$ws.Cells.Item(216, 15).Value = 5500
$ws.Cells.Item(216, 16).Value = 5500
$ws.Cells.Item(216, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(216, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(216, 19).Value = 306
$ws.Cells.Item(216, 20).Value = 18
